$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.8255223037831734
$ws.Range("H2").Value = 0.9466164436939862

$ws.Range("C3").Value = 0.8255223037831734
$ws.Range("H3").Value = 0.8230378317334839

$ws.Range("C4").Value = 0.8255223037831734
$ws.Range("H4").Value = 0.8805122628971849
